# Extend the table from columns A:O to A:Q (add columns P and Q),
# continue the row-1 header sequence, swap a handful of 1/2 values in
# columns I/K/M/O, and fill the two new columns P/Q with 2s for every
# data row (2-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): add P1=14, Q1=15, carrying the s="1" header style ---
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14

$ws.Range("O1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: swap I/K and M/O values, then fill P/Q with 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2  # P (new)
    $ws.Cells.Item($r, 17).Value = 2  # Q (new)
}
